# Weekly update: insert 5 new price rows at the top of the data block
# (right after the existing row 68), shifting all subsequent rows down by 5.
# This mirrors the "Fruta / hortaliza, semanal" weekly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before row 69; Excel copies formatting (incl. the
# date-style on column D) from the row above automatically.
$ws.Range("A69:A73").EntireRow.Insert()

# Common columns shared by every data row in this sheet.
$mercadoId   = 1
$mercado     = "Agrícola del Norte S.A. de Arica"
$region      = "Arica y Parinacota"
$codreg      = 15
$tipo        = "Fruta"
$productoId  = 100109
$producto    = "Uva"
$categoriaId = 100109001
$categoria   = "Uva"
$fecha       = 44994

function Set-DataRow {
    param($Row, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad)

    $ws.Cells.Item($Row, 1).Value2  = $mercadoId
    $ws.Cells.Item($Row, 2).Value2  = $mercado
    $ws.Cells.Item($Row, 3).Value2  = $region
    $ws.Cells.Item($Row, 4).Value2  = $fecha
    $ws.Cells.Item($Row, 5).Value2  = $codreg
    $ws.Cells.Item($Row, 6).Value2  = $tipo
    $ws.Cells.Item($Row, 7).Value2  = $productoId
    $ws.Cells.Item($Row, 8).Value2  = $producto
    $ws.Cells.Item($Row, 9).Value2  = $categoriaId
    $ws.Cells.Item($Row, 10).Value2 = $categoria
    $ws.Cells.Item($Row, 11).Value2 = $Variedad
    $ws.Cells.Item($Row, 12).Value2 = $Calidad
    $ws.Cells.Item($Row, 13).Value2 = $Volumen
    $ws.Cells.Item($Row, 14).Value2 = $PrecioMin
    $ws.Cells.Item($Row, 15).Value2 = $PrecioMax
    $ws.Cells.Item($Row, 16).Value2 = $PrecioProm
    $ws.Cells.Item($Row, 17).Value2 = $Unidad
    $ws.Cells.Item($Row, 18).Value2 = $Origen
    $ws.Cells.Item($Row, 19).Value2 = $PrecioKg
    $ws.Cells.Item($Row, 20).Value2 = $KgUnidad
}

Set-DataRow 69 "Red Globe"         "Primera" 200 17000 18000 17500 "$/bandeja 18 kilos" "Región de O'Higgins"                  972 18
Set-DataRow 70 "Red Globe"         "Segunda" 250 19000 20000 19400 "$/caja 20 kilos"    "Región de O'Higgins"                  970 20
Set-DataRow 71 "Rosada pastilla"   "Primera" 250 17000 18000 17400 "$/bandeja 18 kilos" "Provincia de San Felipe de Aconcagua" 967 18
Set-DataRow 72 "Superior Seedless" "Segunda" 250 19000 20000 19400 "$/caja 20 kilos"    "Región de O'Higgins"                  970 20
Set-DataRow 73 "Superior Seedless" "Tercera" 200 17000 18000 17500 "$/bandeja 18 kilos" "Región de O'Higgins"                  972 18
